# Reflex agent with internal state
# - the "environment" grid (A1:D9) tracks each location's cleanliness; cells
#   that were previously blank (no percept recorded yet) now hold the agent's
#   remembered state for that square: "Clean".
# - selection/view moved to reflect where the agent last looked.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "Clean"
$ws.Range("D5").Value = "Clean"
$ws.Range("C6").Value = "Clean"
$ws.Range("C7").Value = "Clean"
$ws.Range("C8").Value = "Clean"
$ws.Range("D8").Value = "Clean"
$ws.Range("C9").Value = "Clean"
$ws.Range("D9").Value = "Clean"

[void]$ws.Range("G11").Select()
